$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 3 through 6 entirely (delete the extra "proprietaire" rows)
$ws.Range("A3:K6").Clear()

# Reset row 2 (the duplicated/new contrat row) - text columns blank, numeric columns recalculated to 0
$ws.Range("A2:H2").Value = " "
$ws.Range("I2:K2").Value = 0
